$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 47
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# Row 132
$ws.Range("H132").Value = 2988.0557
$ws.Range("I132").Value = 2910.3264
$ws.Range("K132").Value = 8730.9792
$ws.Range("M132").Value = -6200.9792

# Row 135
$ws.Range("H135").Value = 788.46295
$ws.Range("I135").Value = 436.83673
$ws.Range("J135").Value = 4234.4
$ws.Range("K135").Value = 3931.53057
$ws.Range("L135").Value = 38109.6
$ws.Range("M135").Value = -1396.53057
$ws.Range("N135").Value = -43179.6

# Row 137
$ws.Range("H137").Value = 1151.6666
$ws.Range("I137").Value = 969.2826
$ws.Range("J137").Value = 1750.9286
$ws.Range("K137").Value = 2907.8478
$ws.Range("L137").Value = 5252.7858
$ws.Range("M137").Value = -357.8478
$ws.Range("N137").Value = -10352.7858

# Row 141
$ws.Range("H141").Value = 4599.8965
$ws.Range("I141").Value = 1710.3684
$ws.Range("J141").Value = 10090
$ws.Range("K141").Value = 5131.1052
$ws.Range("L141").Value = 30270
$ws.Range("M141").Value = 48.89480000000003
$ws.Range("N141").Value = -40630

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 692881.5600000001
$ws.Range("I32").Value = 795889.8
$ws.Range("J32").Value = 14743.917
$ws.Range("K32").Value = 795889.8
$ws.Range("L32").Value = 14743.917
$ws.Range("M32").Value = -795602.8
$ws.Range("N32").Value = -15317.917

# Row 132
$ws.Range("H132").Value = 2226.7124
$ws.Range("I132").Value = 1557.2909
$ws.Range("J132").Value = 4272.1665
$ws.Range("K132").Value = 4671.8727
$ws.Range("L132").Value = 12816.4995
$ws.Range("M132").Value = -2141.8727
$ws.Range("N132").Value = -17876.4995

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1526.4667
$ws.Range("I20").Value = 1127.6666
$ws.Range("J20").Value = 2124.6667
$ws.Range("K20").Value = 1127.6666
$ws.Range("L20").Value = 2124.6667
$ws.Range("M20").Value = -880.6666
$ws.Range("N20").Value = -2618.6667

# Row 64
$ws.Range("H64").Value = 636.3
$ws.Range("J64").Value = 601.4
$ws.Range("L64").Value = 601.4
$ws.Range("N64").Value = -1051.4

# Row 67
$ws.Range("H67").Value = 636.3
$ws.Range("J67").Value = 601.4
$ws.Range("L67").Value = 601.4
$ws.Range("N67").Value = -2161.4

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 67252.5
$ws.Range("J2").Value = 67252.5
$ws.Range("L2").Value = 67252.5
$ws.Range("N2").Value = -67478.5

# Row 11
$ws.Range("H11").Value = 2598.6667
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 2598.6667
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 2598.6667
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -2878.6667

# Row 41
$ws.Range("H41").Value = 7000
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 9000
$ws.Range("K41").Value = 1000
$ws.Range("L41").Value = 9000
$ws.Range("M41").Value = -572
$ws.Range("N41").Value = -9856

# Row 69
$ws.Range("H69").Value = 30000
$ws.Range("I69").Value = 30000
$ws.Range("J69").Value = 30000
$ws.Range("K69").Value = 30000
$ws.Range("L69").Value = 30000
$ws.Range("M69").Value = -29251
$ws.Range("N69").Value = -31498

# Row 72
$ws.Range("H72").Value = 30000
$ws.Range("I72").Value = 30000
$ws.Range("J72").Value = 30000
$ws.Range("K72").Value = 90000
$ws.Range("L72").Value = 90000
$ws.Range("M72").Value = -86256
$ws.Range("N72").Value = -97488

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 55.25926
$ws.Range("I12").Value = 37.714287
$ws.Range("J12").Value = 61.4
$ws.Range("K12").Value = 113.142861
$ws.Range("L12").Value = 184.2
$ws.Range("M12").Value = 59.857139
$ws.Range("N12").Value = -530.2

# Row 56
$ws.Range("H56").Value = 4753.75
$ws.Range("I56").Value = 4753.75
$ws.Range("K56").Value = 4753.75
$ws.Range("M56").Value = -4223.75

# Row 64
$ws.Range("H64").Value = 2064
$ws.Range("J64").Value = 2071.4285
$ws.Range("L64").Value = 6214.2855
$ws.Range("N64").Value = -6754.2855

# Row 67
$ws.Range("H67").Value = 2064
$ws.Range("J67").Value = 2071.4285
$ws.Range("L67").Value = 6214.2855
$ws.Range("N67").Value = -8086.2855

# Row 80
$ws.Range("H80").Value = 4308
$ws.Range("I80").Value = 2601
$ws.Range("J80").Value = 4429.9287
$ws.Range("K80").Value = 7803
$ws.Range("L80").Value = 13289.7861
$ws.Range("M80").Value = -6867
$ws.Range("N80").Value = -15161.7861

# Row 83
$ws.Range("H83").Value = 4308
$ws.Range("I83").Value = 2601
$ws.Range("J83").Value = 4429.9287
$ws.Range("K83").Value = 23409
$ws.Range("L83").Value = 39869.35830000001
$ws.Range("M83").Value = -18729
$ws.Range("N83").Value = -49229.35830000001

# Row 107
$ws.Range("H107").Value = 16666896
$ws.Range("I107").Value = 283.5625
$ws.Range("J107").Value = 35714452
$ws.Range("K107").Value = 850.6875
$ws.Range("L107").Value = 107143356
$ws.Range("M107").Value = 1069.3125
$ws.Range("N107").Value = -107147196

# Row 113
$ws.Range("H113").Value = 482.0926
$ws.Range("I113").Value = 457.66666
$ws.Range("J113").Value = 497.63635
$ws.Range("K113").Value = 1372.99998
$ws.Range("L113").Value = 1492.90905
$ws.Range("M113").Value = 797.0000199999999
$ws.Range("N113").Value = -5832.90905

# Row 131
$ws.Range("H131").Value = 5526.6924
$ws.Range("I131").Value = 453.75
$ws.Range("J131").Value = 7781.3335
$ws.Range("K131").Value = 1361.25
$ws.Range("L131").Value = 23344.0005
$ws.Range("M131").Value = 3678.75
$ws.Range("N131").Value = -33424.00049999999

$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

# Row 116
$ws.Range("H116").Value = 17999.666
$ws.Range("J116").Value = 17999.666
$ws.Range("L116").Value = 17999.666
$ws.Range("N116").Value = -27177.666

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5322.143
$ws.Range("I7").Value = 4833.8887
$ws.Range("J7").Value = 6201
$ws.Range("K7").Value = 4833.8887
$ws.Range("L7").Value = 6201
$ws.Range("M7").Value = -4721.8887
$ws.Range("N7").Value = -6425

# Row 68
$ws.Range("H68").Value = 1489.2051
$ws.Range("I68").Value = 1446.4412
$ws.Range("K68").Value = 1446.4412
$ws.Range("M68").Value = -697.4412

# Row 71
$ws.Range("H71").Value = 1489.2051
$ws.Range("I71").Value = 1446.4412
$ws.Range("K71").Value = 7232.206
$ws.Range("M71").Value = -3488.206

# Row 106
$ws.Range("H106").Value = 77185
$ws.Range("J106").Value = 77185
$ws.Range("L106").Value = 77185
$ws.Range("N106").Value = -79709

# Row 126
$ws.Range("H126").Value = 5322.143
$ws.Range("I126").Value = 4833.8887
$ws.Range("J126").Value = 6201
$ws.Range("K126").Value = 14501.6661
$ws.Range("L126").Value = 18603
$ws.Range("M126").Value = -12031.6661
$ws.Range("N126").Value = -23543

# Row 132
$ws.Range("H132").Value = 2543.9155
$ws.Range("I132").Value = 2249.875
$ws.Range("J132").Value = 3641.6667
$ws.Range("K132").Value = 6749.625
$ws.Range("L132").Value = 10925.0001
$ws.Range("M132").Value = -4219.625
$ws.Range("N132").Value = -15985.0001

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 3206085.2
$ws.Range("I132").Value = 817.64105
$ws.Range("J132").Value = 12821888
$ws.Range("K132").Value = 2452.92315
$ws.Range("L132").Value = 38465664
$ws.Range("M132").Value = 77.07685000000038
$ws.Range("N132").Value = -38470724
